# Generate Report for handback
# Adds a new handback-status row for file 60eff7f1-41b5-421a-a2b6-54713c617a0f
# to the Overview / zh-cn / de-de worksheets, mirroring the layout already
# used by the existing 4bf8b770... and 02aa88c5... rows.

$wb = $excel.ActiveWorkbook

$uuid     = "60eff7f1-41b5-421a-a2b6-54713c617a0f"
$mdName   = "$uuid.md"
$fileHash = "1479dbe8b70207e844a6c4f3209c6d2e5d9e8c9d"
$zhXlf    = "$uuid.$fileHash.zh-cn.xlf"
$deXlf    = "$uuid.$fileHash.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"
$include      = "Include"

$zhHandoffAt  = "2016-02-19 05:37:11"
$zhHandbackAt = "2016-02-19 05:37:52"
$deHandoffAt  = "2016-02-19 05:37:20"
$deHandbackAt = "2016-02-19 05:38:12"

$srcCommit   = "e4a6a0a4f6b6f5c6d6c3e8c1b6d9a4f3a6b5c6d7"
$zhOffCommit = "1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b"
$zhMdCommit  = "2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c"
$zhBackCommit= "3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c2d"
$deOffCommit = "4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c2d3e"
$deMdCommit  = "5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c2d3e4f"
$deBackCommit= "6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c2d3e4f5a"

$srcUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$mdName"

$zhOffUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhOffCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/$zhXlf"
$zhMdUrl   = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/$zhMdCommit/e2e/$mdName"
$zhBackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$zhBackCommit/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/$zhXlf"

$deOffUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deOffCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/$deXlf"
$deMdUrl   = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/$deMdCommit/e2e/$mdName"
$deBackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$deBackCommit/ol-handback/OpenLocalizationTestOrg/oltest.de-de/terryjin/$deXlf"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $srcUrl, "", "", $mdName)
$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Correspond Handoff File |
#   Correspond Handoff Datetime | Target File | Correspond Handback File |
#   Correspond Handback DateTime | Handoff Reason | Dependency From
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $srcUrl, "", "", $mdName)
$wsZh.Range("B4").Value = $statusInSync
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), $zhOffUrl, "", "", $zhXlf)
$wsZh.Range("D4").Value = $zhHandoffAt
$wsZh.Hyperlinks.Add($wsZh.Range("E4"), $zhMdUrl, "", "", $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), $zhBackUrl, "", "", $zhXlf)
$wsZh.Range("G4").Value = $zhHandbackAt
$wsZh.Range("H4").Value = $include

# ---------------------------------------------------------------------
# Sheet "de-de": same layout as zh-cn
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $srcUrl, "", "", $mdName)
$wsDe.Range("B4").Value = $statusInSync
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), $deOffUrl, "", "", $deXlf)
$wsDe.Range("D4").Value = $deHandoffAt
$wsDe.Hyperlinks.Add($wsDe.Range("E4"), $deMdUrl, "", "", $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), $deBackUrl, "", "", $deXlf)
$wsDe.Range("G4").Value = $deHandbackAt
$wsDe.Range("H4").Value = $include

Write-Host "Added handback row for $uuid to Overview, zh-cn, de-de sheets."
